$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns: Row ID -> Id, Task -> Task Name, Start Date -> Start, End Date -> Finish
$ws.Range("A1").Value = "Id"
$ws.Range("C1").Value = "Task Name"
$ws.Range("E1").Value = "Start"
$ws.Range("F1").Value = "Finish"

# Move the active cell selection from F1 to F2
$ws.Range("F2").Select()
